# "Reviewterminen von WiMis eingetragen"
# Fill in the second-reviewer (and, on HUE, all reviewer) columns with the
# WiMis' initials instead of their bare last names, and leave the HUE sheet
# selected/active as the last-touched sheet.

$wb = $excel.ActiveWorkbook

# --- PUE sheet: column L ("Review 2 (WIMi)") gets the reviewer's initials ---
$pue = $wb.Worksheets.Item("PUE")

$pueL = @{
    3  = "AS"
    4  = "DP"
    5  = "DP"
    6  = "DP"
    7  = "DP"
    8  = "DP"
    9  = "MF"
    10 = "MF"
    11 = "MF"
    12 = "MF"
    13 = "AS"
    14 = "AS"
    15 = "AS"
    16 = "AS"
}

foreach ($row in $pueL.Keys) {
    $pue.Range("L$row").Value = $pueL[$row]
}

# Row 16 previously carried the plain (unbolded) style; match the bold/black
# font already used by rows 13-15 once it holds the same reviewer.
$pue.Range("L16").Font.Color = $pue.Range("L13").Font.Color

# --- HUE sheet: columns J ("Entwurf"), K ("Review 1") and L ("Review 2") ---
$hue = $wb.Worksheets.Item("HUE")

$hueJKL = @{
    3  = @("AS", "MF", "DP")
    4  = @("DP", "AS", "MF")
    5  = @("DP", "AS", "MF")
    6  = @("DP", "AS", "MF")
    7  = @("DP", "AS", "MF")
    8  = @("DP", "AS", "MF")
    9  = @("MF", "DP", "AS")
    10 = @("MF", "DP", "AS")
    11 = @("MF", "DP", "AS")
    12 = @("MF", "DP", "AS")
    13 = @("AS", "MF", "DP")
    14 = @("AS", "MF", "DP")
    15 = @("AS", "MF", "DP")
}

foreach ($row in $hueJKL.Keys) {
    $vals = $hueJKL[$row]
    $hue.Range("J$row").Value = $vals[0]
    $hue.Range("K$row").Value = $vals[1]
    $hue.Range("L$row").Value = $vals[2]
}

# --- Selections / active sheet ---
$null = $pue.Range("L19").Select()

$null = $hue.Activate()
$null = $hue.Range("G21").Select()
